# cambios en los datos
# Remove the two oldest accounts (ES32 4434... / ALTA, ES32 1226... / BAJA)
# and keep the remaining three, renumbering cliente_id and adding the new
# "clasificacion" column values (POOLED / SEGREGADA).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the three surviving accounts (previously rows 4-6) up into rows 2-4,
# dropping the "fecha_cierre" values (none of the survivors had one) and
# filling in the new "clasificacion" column.
$ws.Range("A2").Value = "ES32 5829 354  1000"
$ws.Range("B2").Value = "ALTA"
$ws.Range("C2").Value = 39588
$ws.Range("D2").Clear()
$ws.Range("E2").Value = "POOLED"
$ws.Range("F2").Value = 1

$ws.Range("A3").Value = "ES32 6784 345 0000"
$ws.Range("B3").Value = "ALTA"
$ws.Range("C3").Value = 43785
$ws.Range("D3").Clear()
$ws.Range("E3").Value = "SEGREGADA"
$ws.Range("F3").Value = 2

$ws.Range("A4").Value = "ES32 893 455 2333"
$ws.Range("B4").Value = "ALTA"
$ws.Range("C4").Value = 41311
$ws.Range("D4").Clear()
$ws.Range("E4").Value = "SEGREGADA"
$ws.Range("F4").Value = 3

# Clear out the now-empty old rows 5 and 6 (data columns only - leave the
# style-only anchor cell I6 untouched).
$ws.Range("A5:F6").Clear()

# Update the active selection to match the saved state
$ws.Range("B7").Select()
